$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column (becomes the "TabName" column)
$ws.Columns("A:A").Insert()

# ---- Row 1: headers ----
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# ---- Row 2: CasesTab ----
$casesTabQuery = @'
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
    WHERE a.pubmed_id IN ['31765263'] 
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
'@
$statQueryText = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
  WHERE a.pubmed_id IN ['31765263']
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
'@
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $casesTabQuery
$ws.Range("C2").Value = $statQueryText
$ws.Range("D2").Value = "TC02_Trials_Filter_PubmedID-317_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC02_Trials_Filter_PubmedID-317_WebData.xlsx"
$ws.Range("B2").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Rows("2").RowHeight = 195

# ---- Row 3: FilesTab ----
$filesTabQuery = @'
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
WHERE a.pubmed_id IN ['31765263']
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS `File Name`,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS `File Format`,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS `Trial Code`,
    a.arm_id AS Arm,
    c.case_id AS `Case ID`
'@
$ws.Range("A3").Value = "FilesTab"
$ws.Range("B3").Value = $filesTabQuery
$ws.Range("C3").Value = $statQueryText
$ws.Range("D3").Value = "TC02_Trials_Filter_PubmedID-317_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC02_Trials_Filter_PubmedID-317_WebData.xlsx"
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true
$ws.Rows("3").RowHeight = 409.5

# ---- Column widths ----
$ws.Columns("A:A").ColumnWidth = 8
$ws.Columns("B:C").ColumnWidth = 74.95
$ws.Columns("D:D").ColumnWidth = 69.5
$ws.Columns("E:E").ColumnWidth = 27.6

# ---- Selection ----
$null = $ws.Range("C3").Select()
